$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("7520950", "Can - Pizza Sauce", "2", "41.22", "82.44"),
    @("2742112", "Relish - Sweet", "1", "34.58", "34.58"),
    @("8255796", "Tuna White Chunk (Pouch)", "5", "72.00", "360.00"),
    @("1030192", "BBQ - Sauce", "1", "72.94", "72.94"),
    @("3427135", "Dolmades (Stuffed Grape Leaves)", "2", "54.02", "108.04"),
    @("1028165", "Pickle - Dill Chip", "1", "34.09", "34.09"),
    @("1035842", "Feta - Pail", "1", "92.87", "92.87"),
    @("2212975", "Roast Beef (Unsliced)", "6", "119.20", "715.20"),
    @("6364494", "Yogurt - Greek (Bulk)", "3", "26.96", "80.88"),
    @("1365278", "Vegan Chicken Tenders", "1", "87.80", "87.80"),
    @("1028188", "Tortellini - Cheese", "3", "32.22", "96.66"),
    @("6318732", "FRZ Fruit - Apple (Sliced)", "1", "66.72", "66.72"),
    @("2825368", "Sausage - Chicken Patty", "3", "50.98", "152.94"),
    @("9546982", "Arugula - Fresh", "5", "20.35", "101.75"),
    @("7489339", "Garlic - Fresh (Peeled)", "1", "70.08", "70.08"),
    @("89213", "Potato - Chef", "1", "17.47", "17.47"),
    @("6264014", "Tomato - Grape", "2", "18.50", "37.00")
)

$startRow = 12
$endRow = $startRow + $newRows.Count - 1

# Make sure these new cells are treated as plain text (matching the
# existing rows, which store SKU/Quantity/Cost/Total as text strings
# rather than numbers), then fill in the values.
$ws.Range("A$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
